$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting existing data to the right
$ws.Columns("A").Insert()

# New header for the inserted column
$ws.Range("A1").Value = "model_id"

# New model_id values for rows 2-7
$ws.Range("A2").Value = 10
$ws.Range("A3").Value = 20
$ws.Range("A4").Value = 30
$ws.Range("A5").Value = 40
$ws.Range("A6").Value = 50
$ws.Range("A7").Value = 60

# Update the active selection to match the target state
$ws.Range("A8").Select()
